$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.97977521793026
$ws.Range("C2").Value = 4.577166254778756
$ws.Range("D2").Value = 8.869783142047433
$ws.Range("E2").Value = 13.5665949532071
$ws.Range("F2").Value = 34.7032031195418
$ws.Range("J2").Value = 10.01107948678177
$ws.Range("K2").Value = 10.32542902512413
$ws.Range("O2").Value = 26.3942880787376

$ws.Range("B3").Value = 10.70408834924434
$ws.Range("C3").Value = 4.368798852750527
$ws.Range("D3").Value = 8.802820155991634
$ws.Range("E3").Value = 13.50085619410533
$ws.Range("F3").Value = 34.76823429456962
$ws.Range("J3").Value = 10.0182502811835
$ws.Range("K3").Value = 10.13360073497342
$ws.Range("O3").Value = 26.4837929842765

$ws.Range("B4").Value = 10.53302183520635
$ws.Range("C4").Value = 4.234880977079799
$ws.Range("D4").Value = 8.762985515714593
$ws.Range("E4").Value = 13.4631686560127
$ws.Range("F4").Value = 34.81645626696238
$ws.Range("J4").Value = 10.02423992581137
$ws.Range("K4").Value = 10.01547233748019
$ws.Range("O4").Value = 26.54431504987409

$ws.Range("B5").Value = 10.46296495252434
$ws.Range("C5").Value = 4.178850349714833
$ws.Range("D5").Value = 8.74708901222974
$ws.Range("E5").Value = 13.44849594336944
$ws.Range("F5").Value = 34.83818776045587
$ws.Range("J5").Value = 10.02708004871977
$ws.Range("K5").Value = 9.967309631941427
$ws.Range("O5").Value = 26.57037491010158

$ws.Range("B6").Value = 10.45131433700268
$ws.Range("C6").Value = 4.169459925273874
$ws.Range("D6").Value = 8.74447013786335
$ws.Range("E6").Value = 13.44610127941544
$ws.Range("F6").Value = 34.84192177841585
$ws.Range("J6").Value = 10.0275757718844
$ws.Range("K6").Value = 9.959312621812147
$ws.Range("O6").Value = 26.57478639622617

$ws.Range("B7").Value = 10.53207828158484
$ws.Range("C7").Value = 4.234131166191173
$ws.Range("D7").Value = 8.762769749157915
$ws.Range("E7").Value = 13.46296798453069
$ws.Range("F7").Value = 34.81674092756337
$ws.Range("J7").Value = 10.0242766116644
$ws.Range("K7").Value = 10.01482281401891
$ws.Range("O7").Value = 26.54466085130588

$ws.Range("B8").Value = 10.88515218942337
$ws.Range("C8").Value = 4.506587323250848
$ws.Range("D8").Value = 8.846437222846431
$ws.Range("E8").Value = 13.5433804046697
$ws.Range("F8").Value = 34.72390216692591
$ws.Range("J8").Value = 10.01322285175264
$ws.Range("K8").Value = 10.25939523624328
$ws.Range("O8").Value = 26.42399248533813

$ws.Range("B9").Value = 11.55894305067962
$ws.Range("C9").Value = 4.991776519014489
$ws.Range("D9").Value = 9.02001783931829
$ws.Range("E9").Value = 13.72172295503272
$ws.Range("F9").Value = 34.60782430427126
$ws.Range("J9").Value = 10.00412055979892
$ws.Range("K9").Value = 10.73359960016869
$ws.Range("O9").Value = 26.23165078213475

$ws.Range("B10").Value = 12.03722541322798
$ws.Range("C10").Value = 5.316523339522027
$ws.Range("D10").Value = 9.152438314685686
$ws.Range("E10").Value = 13.86450710326222
$ws.Range("F10").Value = 34.56297518525514
$ws.Range("J10").Value = 10.00507235230126
$ws.Range("K10").Value = 11.07527318183077
$ws.Range("O10").Value = 26.11750666208516

$ws.Range("B11").Value = 12.25014878960479
$ws.Range("C11").Value = 5.457093768452164
$ws.Range("D11").Value = 9.213542207880879
$ws.Range("E11").Value = 13.93182211561508
$ws.Range("F11").Value = 34.55138281419045
$ws.Range("J11").Value = 10.00715630411976
$ws.Range("K11").Value = 11.2285694728995
$ws.Range("O11").Value = 26.07151366382426

$ws.Range("B12").Value = 12.33003226294401
$ws.Range("C12").Value = 5.509276886782436
$ws.Range("D12").Value = 9.236787980637578
$ws.Range("E12").Value = 13.95763483059312
$ws.Range("F12").Value = 34.54826113960534
$ws.Range("J12").Value = 10.00818196254253
$ws.Range("K12").Value = 11.28625910274023
$ws.Range("O12").Value = 26.0549529921021

$ws.Range("B13").Value = 12.3128622689701
$ws.Range("C13").Value = 5.498085179962215
$ws.Range("D13").Value = 9.231777103428254
$ws.Range("E13").Value = 13.95206156713824
$ws.Range("F13").Value = 34.54887703675679
$ws.Range("J13").Value = 10.00795056534113
$ws.Range("K13").Value = 11.27385145784371
$ws.Range("O13").Value = 26.05848151872897

$ws.Range("B14").Value = 12.25673622195429
$ws.Range("C14").Value = 5.461407999152593
$ws.Range("D14").Value = 9.215452605751672
$ws.Range("E14").Value = 13.93393939324786
$ws.Range("F14").Value = 34.55110057289701
$ws.Range("J14").Value = 10.00723595119376
$ws.Range("K14").Value = 11.23332315209587
$ws.Range("O14").Value = 26.07013403812778

$ws.Range("B15").Value = 12.22225806195776
$ws.Range("C15").Value = 5.438805216172656
$ws.Range("D15").Value = 9.205466804740606
$ws.Range("E15").Value = 13.92288043680712
$ws.Range("F15").Value = 34.55262772224332
$ws.Range("J15").Value = 10.00682900141034
$ws.Range("K15").Value = 11.20844992506682
$ws.Range("O15").Value = 26.07738308741204

$ws.Range("B16").Value = 12.0232100447516
$ws.Range("C16").Value = 5.307191025265074
$ws.Range("D16").Value = 9.14846095946986
$ws.Range("E16").Value = 13.86015393903859
$ws.Range("F16").Value = 34.56391017506321
$ws.Range("J16").Value = 10.00496930909481
$ws.Range("K16").Value = 11.06520733722165
$ws.Range("O16").Value = 26.12063202522345

$ws.Range("B17").Value = 11.89985253210395
$ws.Range("C17").Value = 5.224602246393103
$ws.Range("D17").Value = 9.113698993601888
$ws.Range("E17").Value = 13.82226580654793
$ws.Range("F17").Value = 34.57308913077918
$ws.Range("J17").Value = 10.00425070300985
$ws.Range("K17").Value = 10.9767472815146
$ws.Range("O17").Value = 26.14868521062034

$ws.Range("B18").Value = 11.82846711632168
$ws.Range("C18").Value = 5.176426432950647
$ws.Range("D18").Value = 9.093787310146258
$ws.Range("E18").Value = 13.80069695207794
$ws.Range("F18").Value = 34.57919780809888
$ws.Range("J18").Value = 10.00399283892407
$ws.Range("K18").Value = 10.925669679058
$ws.Range("O18").Value = 26.16537881760951

$ws.Range("B19").Value = 11.80422526160321
$ws.Range("C19").Value = 5.159999973854727
$ws.Range("D19").Value = 9.087060249333515
$ws.Range("E19").Value = 13.79343300536564
$ws.Range("F19").Value = 34.58140846291414
$ws.Range("J19").Value = 10.00393225617855
$ws.Range("K19").Value = 10.9083434060719
$ws.Range("O19").Value = 26.17112675967829

$ws.Range("B20").Value = 11.91302959003661
$ws.Range("C20").Value = 5.233463753594974
$ws.Range("D20").Value = 9.117391053266129
$ws.Range("E20").Value = 13.82627606815951
$ws.Range("F20").Value = 34.57202618599906
$ws.Range("J20").Value = 10.0043111165322
$ws.Range("K20").Value = 10.9861848667382
$ws.Range("O20").Value = 26.14564111258011

$ws.Range("B21").Value = 12.273242632958
$ws.Range("C21").Value = 5.472209544859568
$ws.Range("D21").Value = 9.220244739705581
$ws.Range("E21").Value = 13.93925371987459
$ws.Range("F21").Value = 34.5504130441273
$ws.Range("J21").Value = 10.00743943960374
$ws.Range("K21").Value = 11.24523749045782
$ws.Range("O21").Value = 26.06668815989954

$ws.Range("B22").Value = 12.50428313048737
$ws.Range("C22").Value = 5.622129781166121
$ws.Range("D22").Value = 9.288081475575465
$ws.Range("E22").Value = 14.01495934151377
$ws.Range("F22").Value = 34.54367912882082
$ws.Range("J22").Value = 10.01086204487732
$ws.Range("K22").Value = 11.41242021566395
$ws.Range("O22").Value = 26.02007781464031

$ws.Range("B23").Value = 12.38139709567186
$ws.Range("C23").Value = 5.542679194560701
$ws.Range("D23").Value = 9.251825127710413
$ws.Range("E23").Value = 13.97438885320457
$ws.Range("F23").Value = 34.54659660205261
$ws.Range("J23").Value = 10.00890957091057
$ws.Range("K23").Value = 11.32340286451783
$ws.Range("O23").Value = 26.04449712345213

$ws.Range("B24").Value = 11.90707368624258
$ws.Range("C24").Value = 5.229459627760388
$ws.Range("D24").Value = 9.11572164304704
$ws.Range("E24").Value = 13.82446236242051
$ws.Range("F24").Value = 34.57250415311155
$ws.Range("J24").Value = 10.00428331985608
$ws.Range("K24").Value = 10.98191881789159
$ws.Range("O24").Value = 26.14701558812476

$ws.Range("B25").Value = 11.37923818799941
$ws.Range("C25").Value = 4.865981748179183
$ws.Range("D25").Value = 8.972133776867185
$ws.Range("E25").Value = 13.67134847090463
$ws.Range("F25").Value = 34.6321384361872
$ws.Range("J25").Value = 10.00523908326169
$ws.Range("K25").Value = 10.6062560605794
$ws.Range("O25").Value = 26.27892499449896
